$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 148: 四方坪站充电量(kw)
$ws.Range("A148").Value = 45974
$ws.Range("B148").Value = "四方坪站充电量(kw)"
$ws.Range("C148").Value = 600.55099999999993
$ws.Range("D148").Value = 1166.8229999999999
$ws.Range("E148").Value = 300.40199999999999
$ws.Range("F148").Value = 214.78999999999996
$ws.Range("G148").Value = 238.92
$ws.Range("H148").Value = 563.19399999999996
$ws.Range("I148").Value = 490.38300000000004
$ws.Range("J148").Value = 194.34100000000001
$ws.Range("K148").Value = 103.79999999999998
$ws.Range("L148").Value = 244.18799999999999
$ws.Range("M148").Value = 268.65999999999997
$ws.Range("N148").Value = 415.89800000000002
$ws.Range("O148").Value = 809.74900000000002
$ws.Range("P148").Value = 1317.7050000000002
$ws.Range("Q148").Value = 290.13400000000001
$ws.Range("R148").Value = 288.78700000000003
$ws.Range("S148").Value = 408.07799999999992
$ws.Range("T148").Value = 189.14
$ws.Range("U148").Value = 209.76
$ws.Range("V148").Value = 74.11
$ws.Range("W148").Value = 0
$ws.Range("X148").Value = 21.1
$ws.Range("Y148").Value = 26.28
$ws.Range("Z148").Value = 103.57000000000001

# Row 149: 高岭站充电量(kw)
$ws.Range("A149").Value = 45974
$ws.Range("B149").Value = "高岭站充电量(kw)"
$ws.Range("C149").Value = 383.89600000000002
$ws.Range("D149").Value = 380.11199999999997
$ws.Range("E149").Value = 105.59399999999999
$ws.Range("F149").Value = 0
$ws.Range("G149").Value = 32.75
$ws.Range("H149").Value = 67.085999999999999
$ws.Range("I149").Value = 326.99699999999996
$ws.Range("J149").Value = 101.76400000000001
$ws.Range("K149").Value = 254.80099999999999
$ws.Range("L149").Value = 212.46800000000002
$ws.Range("M149").Value = 275.62700000000001
$ws.Range("N149").Value = 193.02699999999999
$ws.Range("O149").Value = 364.81700000000001
$ws.Range("P149").Value = 343.14599999999996
$ws.Range("Q149").Value = 354.86799999999994
$ws.Range("R149").Value = 406.72899999999998
$ws.Range("S149").Value = 101.32900000000001
$ws.Range("T149").Value = 109.39
$ws.Range("U149").Value = 49.675000000000004
$ws.Range("V149").Value = 48.998000000000005
$ws.Range("W149").Value = 34.008000000000003
$ws.Range("X149").Value = 38.125
$ws.Range("Y149").Value = 0
$ws.Range("Z149").Value = 11.852

# Match the author's final view state: scrolled/selected at G151
$ws.Range("G151").Select() | Out-Null
